$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new value in D4 (new shared string, no wrap-text style like C4 has)
$ws.Range("D4").Value = '("PSOT","L","B505","RBS","B1")'

# Update the selection to D4 as shown in the diff
$ws.Range("D4").Select()
